# Scheduled runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, Leve price/profit columns H:N) across the ALC, ARM,
# BSM, CRP, CUL, GSM, LTW and WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 837.9286
$ws.Range("I32").Value = 881.5
$ws.Range("J32").Value = 805.25
$ws.Range("K32").Value = 881.5
$ws.Range("L32").Value = 805.25
$ws.Range("M32").Value = -555.5
$ws.Range("N32").Value = -1457.25
$ws.Range("H108").Value = 35383.332
$ws.Range("J108").Value = 35383.332
$ws.Range("L108").Value = 35383.332
$ws.Range("N108").Value = -43063.332
$ws.Range("H110").Value = 44816.5
$ws.Range("J110").Value = 44816.5
$ws.Range("L110").Value = 44816.5
$ws.Range("N110").Value = -52996.5
$ws.Range("H113").Value = 41669460
$ws.Range("I113").Value = 62502124
$ws.Range("J113").Value = 4124.5
$ws.Range("K113").Value = 62502124
$ws.Range("L113").Value = 4124.5
$ws.Range("M113").Value = -62498870
$ws.Range("N113").Value = -10632.5
$ws.Range("H115").Value = 1049.375
$ws.Range("I115").Value = 770.7143
$ws.Range("K115").Value = 2312.1429
$ws.Range("M115").Value = -745.1428999999998
$ws.Range("H116").Value = 4030.4688
$ws.Range("I116").Value = 2914.2273
$ws.Range("J116").Value = 6486.2
$ws.Range("K116").Value = 2914.2273
$ws.Range("L116").Value = 6486.2
$ws.Range("M116").Value = 527.7727
$ws.Range("N116").Value = -13370.2
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H129").Value = 1207.8529
$ws.Range("I129").Value = 604.75
$ws.Range("J129").Value = 1743.9445
$ws.Range("K129").Value = 1814.25
$ws.Range("L129").Value = 5231.833500000001
$ws.Range("M129").Value = 3185.75
$ws.Range("N129").Value = -15231.8335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3425.4814
$ws.Range("I45").Value = 2374.1875
$ws.Range("J45").Value = 4954.636
$ws.Range("K45").Value = 2374.1875
$ws.Range("L45").Value = 4954.636
$ws.Range("M45").Value = -1997.1875
$ws.Range("N45").Value = -5708.636
$ws.Range("H97").Value = 965.10254
$ws.Range("I97").Value = 809.4167
$ws.Range("J97").Value = 2833.3333
$ws.Range("K97").Value = 809.4167
$ws.Range("L97").Value = 2833.3333
$ws.Range("M97").Value = -313.4167
$ws.Range("N97").Value = -3825.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1929.6046
$ws.Range("I86").Value = 1966.6471
$ws.Range("J86").Value = 1789.6666
$ws.Range("K86").Value = 1966.6471
$ws.Range("L86").Value = 1789.6666
$ws.Range("M86").Value = -843.6470999999999
$ws.Range("N86").Value = -4035.6666
$ws.Range("H89").Value = 1929.6046
$ws.Range("I89").Value = 1966.6471
$ws.Range("J89").Value = 1789.6666
$ws.Range("K89").Value = 9833.235499999999
$ws.Range("L89").Value = 8948.333000000001
$ws.Range("M89").Value = -4217.235499999999
$ws.Range("N89").Value = -20180.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3295.3333
$ws.Range("I62").Value = 3239.0908
$ws.Range("J62").Value = 3450
$ws.Range("K62").Value = 3239.0908
$ws.Range("L62").Value = 3450
$ws.Range("M62").Value = -2615.0908
$ws.Range("N62").Value = -4698
$ws.Range("H65").Value = 3295.3333
$ws.Range("I65").Value = 3239.0908
$ws.Range("J65").Value = 3450
$ws.Range("K65").Value = 16195.454
$ws.Range("L65").Value = 17250
$ws.Range("M65").Value = -13075.454
$ws.Range("N65").Value = -23490
$ws.Range("H86").Value = 62510124
$ws.Range("I86").Value = 76934070
$ws.Range("J86").Value = 6336
$ws.Range("K86").Value = 76934070
$ws.Range("L86").Value = 6336
$ws.Range("M86").Value = -76932947
$ws.Range("N86").Value = -8582
$ws.Range("H89").Value = 62510124
$ws.Range("I89").Value = 76934070
$ws.Range("J89").Value = 6336
$ws.Range("K89").Value = 384670350
$ws.Range("L89").Value = 31680
$ws.Range("M89").Value = -384664734
$ws.Range("N89").Value = -42912
$ws.Range("H134").Value = 5477.0884
$ws.Range("I134").Value = 6007.773
$ws.Range("J134").Value = 4504.1665
$ws.Range("K134").Value = 18023.319
$ws.Range("L134").Value = 13512.4995
$ws.Range("M134").Value = -15488.319
$ws.Range("N134").Value = -18582.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 11291.25
$ws.Range("I3").Value = 7682.5
$ws.Range("J3").Value = 14900
$ws.Range("K3").Value = 23047.5
$ws.Range("L3").Value = 44700
$ws.Range("M3").Value = -22935.5
$ws.Range("N3").Value = -44924
$ws.Range("H13").Value = 9009362
$ws.Range("I13").Value = 205.85715
$ws.Range("K13").Value = 617.5714499999999
$ws.Range("M13").Value = -449.5714499999999
$ws.Range("H58").Value = 5884330
$ws.Range("I58").Value = 7354537
$ws.Range("K58").Value = 22063611
$ws.Range("M58").Value = -22063483
$ws.Range("H64").Value = 6381.12
$ws.Range("I64").Value = 910
$ws.Range("J64").Value = 7127.1816
$ws.Range("K64").Value = 2730
$ws.Range("L64").Value = 21381.5448
$ws.Range("M64").Value = -2460
$ws.Range("N64").Value = -21921.5448
$ws.Range("H67").Value = 6381.12
$ws.Range("I67").Value = 910
$ws.Range("J67").Value = 7127.1816
$ws.Range("K67").Value = 2730
$ws.Range("L67").Value = 21381.5448
$ws.Range("M67").Value = -1794
$ws.Range("N67").Value = -23253.5448
$ws.Range("H70").Value = 7596.1665
$ws.Range("J70").Value = 7092.75
$ws.Range("L70").Value = 21278.25
$ws.Range("N70").Value = -21908.25
$ws.Range("H73").Value = 7596.1665
$ws.Range("J73").Value = 7092.75
$ws.Range("L73").Value = 21278.25
$ws.Range("N73").Value = -23462.25
$ws.Range("H76").Value = 7436.364
$ws.Range("J76").Value = 8200
$ws.Range("L76").Value = 24600
$ws.Range("N76").Value = -25366
$ws.Range("H79").Value = 7436.364
$ws.Range("J79").Value = 8200
$ws.Range("L79").Value = 24600
$ws.Range("N79").Value = -27252
$ws.Range("H82").Value = 7783.609
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 7783.609
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 23350.827
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -24162.827
$ws.Range("H85").Value = 7783.609
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 7783.609
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 23350.827
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -26158.827
$ws.Range("H94").Value = 7001.091
$ws.Range("I94").Value = 6012
$ws.Range("J94").Value = 7220.8887
$ws.Range("K94").Value = 18036
$ws.Range("L94").Value = 21662.6661
$ws.Range("M94").Value = -17360
$ws.Range("N94").Value = -23014.6661
$ws.Range("H100").Value = 12828469
$ws.Range("J100").Value = 12828469
$ws.Range("L100").Value = 38485407
$ws.Range("N100").Value = -38487029
$ws.Range("H106").Value = 3985.8
$ws.Range("J106").Value = 3985.8
$ws.Range("L106").Value = 11957.4
$ws.Range("N106").Value = -13849.4
$ws.Range("H109").Value = 4788.4707
$ws.Range("I109").Value = 3234.25
$ws.Range("J109").Value = 6170
$ws.Range("K109").Value = 9702.75
$ws.Range("L109").Value = 18510
$ws.Range("M109").Value = -8662.75
$ws.Range("N109").Value = -20590
$ws.Range("H112").Value = 41718660
$ws.Range("I112").Value = 1200
$ws.Range("J112").Value = 46353932
$ws.Range("K112").Value = 3600
$ws.Range("L112").Value = 139061796
$ws.Range("M112").Value = -2492
$ws.Range("N112").Value = -139064012
$ws.Range("H122").Value = 3242.5557
$ws.Range("I122").Value = 361.875
$ws.Range("J122").Value = 3495.8022
$ws.Range("K122").Value = 3256.875
$ws.Range("L122").Value = 31462.2198
$ws.Range("M122").Value = -806.875
$ws.Range("N122").Value = -36362.2198
$ws.Range("H125").Value = 2647.6191
$ws.Range("J125").Value = 2768.4211
$ws.Range("L125").Value = 8305.263300000001
$ws.Range("N125").Value = -18145.2633
$ws.Range("H133").Value = 2512.5
$ws.Range("I133").Value = 2512.5
$ws.Range("K133").Value = 7537.5
$ws.Range("M133").Value = -2477.5
$ws.Range("H134").Value = 1365.7059
$ws.Range("I134").Value = 1365.7059
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4097.1177
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 972.8823000000002
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 44.444443
$ws.Range("I2").Value = 19.5
$ws.Range("J2").Value = 64.40000000000001
$ws.Range("K2").Value = 19.5
$ws.Range("L2").Value = 64.40000000000001
$ws.Range("M2").Value = 93.5
$ws.Range("N2").Value = -290.4
$ws.Range("H102").Value = 1631.2565
$ws.Range("I102").Value = 1201.1852
$ws.Range("J102").Value = 2598.9167
$ws.Range("K102").Value = 1201.1852
$ws.Range("L102").Value = 2598.9167
$ws.Range("M102").Value = 420.8148000000001
$ws.Range("N102").Value = -5842.9167
$ws.Range("H126").Value = 4281.3057
$ws.Range("I126").Value = 1941.9048
$ws.Range("J126").Value = 7556.467
$ws.Range("K126").Value = 5825.7144
$ws.Range("L126").Value = 22669.401
$ws.Range("M126").Value = -3355.7144
$ws.Range("N126").Value = -27609.401

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1039.5
$ws.Range("I22").Value = 916.9167
$ws.Range("J22").Value = 1131.4375
$ws.Range("K22").Value = 916.9167
$ws.Range("L22").Value = 1131.4375
$ws.Range("M22").Value = -621.9167
$ws.Range("N22").Value = -1721.4375
$ws.Range("H27").Value = 1039.5
$ws.Range("I27").Value = 916.9167
$ws.Range("J27").Value = 1131.4375
$ws.Range("K27").Value = 916.9167
$ws.Range("L27").Value = 1131.4375
$ws.Range("M27").Value = -809.9167
$ws.Range("N27").Value = -1345.4375
$ws.Range("H40").Value = 1466.6957
$ws.Range("I40").Value = 1313.2
$ws.Range("K40").Value = 1313.2
$ws.Range("M40").Value = -1177.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 10023335
$ws.Range("J11").Value = 7535002.5
$ws.Range("L11").Value = 7535002.5
$ws.Range("N11").Value = -7535286.5
$ws.Range("H126").Value = 2580.1052
$ws.Range("I126").Value = 1852
$ws.Range("J126").Value = 3828.2856
$ws.Range("K126").Value = 5556
$ws.Range("L126").Value = 11484.8568
$ws.Range("M126").Value = -3086
$ws.Range("N126").Value = -16424.8568
